# Updated cryptos list values (Price + Volume(1h)) per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (values that cannot be misparsed as numbers).
$ws.Range("D2").Value = "41.862.91"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.208.60"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("E7").Value = "  -2.38%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("E10").Value = "  -3.06%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").Value = "2.537.90"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "2.215.98"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").Value = "41.858.75"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").Value = "0.0₃0925"
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("E26").Value = "  -3.04%  "
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("E32").Value = "  -5.66%  "
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("E37").Value = "  -5.14%  "
$ws.Range("E38").Value = "  -7.83%  "
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("D47").Value = "1.456.12"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("E48").Value = "  -11.84%  "
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("E51").Value = "  -2.27%  "

# Numeric-looking text updates: force text entry (NumberFormat "@") so the
# original decimal formatting (trailing zeros, etc.) is preserved exactly as
# a string, then restore the default "Normal" style so no stray formatting
# is introduced.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.57"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.402"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0895"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.796"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.140"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0647"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000240"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0242"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0954"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "96.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.07"
$ws.Range("D51").Style = "Normal"
